# Handback status report generation update.
# Updates the "Latest HO Xliff Generate Date" on the Overview sheet and the
# "Correspond Handoff Datetime" / "Correspond Handback DateTime" values on the
# per-locale sheets for the row corresponding to
# 02ddfda7-a291-4a8b-852d-4a174f942789.md, reflecting a newly generated
# handback report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-20 22:53:45"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-20 22:53:41"
$zhcn.Range("K2").Value = "2016-08-20 22:53:58"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-20 22:53:45"
$dede.Range("K2").Value = "2016-08-20 22:54:09"
